# Update the "想去人数" (want-to-go count) figures refreshed by the
# gh-pages data generation run (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 13540
$ws1.Range("F5").Value = 13358
$ws1.Range("F6").Value = 1035
$ws1.Range("F7").Value = 789
$ws1.Range("F13").Value = 716
$ws1.Range("F15").Value = 43
$ws1.Range("F21").Value = 324
$ws1.Range("F23").Value = 470
$ws1.Range("F25").Value = 52

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 843

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 13540
$ws4.Range("F7").Value = 13358
$ws4.Range("F8").Value = 1035
$ws4.Range("F9").Value = 789
$ws4.Range("F15").Value = 716
$ws4.Range("F19").Value = 43
$ws4.Range("F28").Value = 324
$ws4.Range("F30").Value = 470
$ws4.Range("F33").Value = 843
$ws4.Range("F36").Value = 52
